$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.527.06'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '2.283.53'
$ws.Range("E3").Value = '  -1.20%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '95.39'
$ws.Range("E5").Value = '  +1.50%  '
$ws.Range("D6").Value = '267.74'
$ws.Range("E6").Value = '  -1.02%  '
$ws.Range("E7").Value = '  -1.06%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '0.607'
$ws.Range("D10").Value = '45.19'
$ws.Range("E10").Value = '  +0.87%  '
$ws.Range("D11").Value = '''0.0930'
$ws.Range("E11").Value = '  -0.88%  '
$ws.Range("D12").Value = '7.89'
$ws.Range("E12").Value = '  -2.91%  '
$ws.Range("E13").Value = '  +1.46%  '
$ws.Range("D14").Value = '2.626.83'
$ws.Range("E14").Value = '  -1.21%  '
$ws.Range("D15").Value = '15.21'
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("D16").Value = '0.846'
$ws.Range("E16").Value = '  -1.38%  '
$ws.Range("D17").Value = '2.279.72'
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("D18").Value = '43.505.39'
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("E19").Value = '  +2.26%  '
$ws.Range("D20").Value = '6.19'
$ws.Range("E20").Value = '  -1.50%  '
$ws.Range("D21").Value = '71.92'
$ws.Range("E21").Value = '  +0.44%  '
$ws.Range("D22").Value = '2.59'
$ws.Range("E22").Value = '  +12.80%  '
$ws.Range("D23").Value = '232.02'
$ws.Range("E23").Value = '  -2.90%  '
$ws.Range("D24").Value = '9.15'
$ws.Range("E24").Value = '  -5.45%  '
$ws.Range("E25").Value = '  -0.13%  '
$ws.Range("D26").Value = '2.57'
$ws.Range("E26").Value = '  +2.27%  '
$ws.Range("D27").Value = '11.17'
$ws.Range("E27").Value = '  -1.58%  '
$ws.Range("D28").Value = '3.47'
$ws.Range("E28").Value = '  +2.38%  '
$ws.Range("D29").Value = '39.98'
$ws.Range("E29").Value = '  +2.50%  '
$ws.Range("D30").Value = '2.22'
$ws.Range("E30").Value = '  -6.33%  '
$ws.Range("D31").Value = '174.97'
$ws.Range("E31").Value = '  +1.81%  '
$ws.Range("D32").Value = '21.75'
$ws.Range("E32").Value = '  -3.61%  '
$ws.Range("D33").Value = '''0.0890'
$ws.Range("E33").Value = '  -0.98%  '
$ws.Range("E34").Value = '  -4.00%  '
$ws.Range("E35").Value = '  -0.83%  '
$ws.Range("D37").Value = '0.0351'
$ws.Range("E37").Value = '  -1.94%  '
$ws.Range("D38").Value = '4.35'
$ws.Range("E38").Value = '  -3.16%  '
$ws.Range("E39").Value = '  -4.51%  '
$ws.Range("E40").Value = '  +1.08%  '
$ws.Range("E41").Value = '  +0.23%  '
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").Value = '65.61'
$ws.Range("E42").Value = '  +6.24%  '
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").Value = '12.26'
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("E44").Value = '  +0.56%  '
$ws.Range("E45").Value = '  -2.55%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '0.102'
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").Value = '5.14'
$ws.Range("E47").Value = '  -5.95%  '
$ws.Range("E48").Value = '  -1.86%  '
$ws.Range("D49").Value = '''96.30'
$ws.Range("E49").Value = '  -4.11%  '
$ws.Range("B50").Value = 'WOONetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D50").Value = '0.431'
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").Value = '0.185'
$ws.Range("E51").Value = '  +7.23%  '
